$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# " Since we originally had 3 big stories we simply added another one as 4
# big stories are fine for one iteration, but we also changed priorities..."
# becomes
# " Since we originally had 3 big stories planned in Iteration2, we simply
# added another one as 4 big stories are fine for one iteration, but we
# changed priorities..."
#
# i.e. insert "planned in Iteration2, " right before "we simply added"
# and drop the word "also " from "we also changed priorities".

$r1 = $d.Content
$r1.Find.Execute(
    "Since we originally had 3 big stories we simply added",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Since we originally had 3 big stories planned in Iteration2, we simply added",
    2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute(
    "we also changed priorities",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "we changed priorities",
    2) | Out-Null

Write-Output "done"
